$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ace"
$ws.Range("C2").Value = "Bdkrb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 32.50384966666667
$ws.Range("H2").Value = 97.511549
$ws.Range("I2").Value = 0.1162664629566559
$ws.Range("J2").Value = 0.116266462956656
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.3006476666666667
$ws.Range("N2").Value = 0.9019430000000001
$ws.Range("O2").Value = 0.2555551336960822
$ws.Range("P2").Value = 0.2555551336960821
$ws.Range("Q2").Value = 9.772206559967445
$ws.Range("R2").Value = 87.94985903970701
$ws.Range("S2").Value = 0.0297124914852588
$ws.Range("T2").Value = 0.0297124914852588

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ace"
$ws.Range("C3").Value = "Bdkrb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 32.50384966666667
$ws.Range("H3").Value = 97.511549
$ws.Range("I3").Value = 0.1162664629566559
$ws.Range("J3").Value = 0.116266462956656
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8000553333333333
$ws.Range("N3").Value = 2.400166
$ws.Range("O3").Value = 0.6800593197383766
$ws.Range("P3").Value = 0.6800593197383766
$ws.Range("Q3").Value = 26.00487827968156
$ws.Range("R3").Value = 234.043904517134
$ws.Range("S3").Value = 0.07906809170669062
$ws.Range("T3").Value = 0.07906809170669062

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ace"
$ws.Range("C4").Value = "Bdkrb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 32.50384966666667
$ws.Range("H4").Value = 97.511549
$ws.Range("I4").Value = 0.1162664629566559
$ws.Range("J4").Value = 0.116266462956656
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.07574633333333333
$ws.Range("N4").Value = 0.227239
$ws.Range("O4").Value = 0.06438554656554128
$ws.Range("P4").Value = 0.06438554656554128
$ws.Range("Q4").Value = 2.462047431467889
$ws.Range("R4").Value = 22.158426883211
$ws.Range("S4").Value = 0.007485879764706552
$ws.Range("T4").Value = 0.007485879764706553

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ace"
$ws.Range("C5").Value = "Bdkrb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 245.8810833333334
$ws.Range("H5").Value = 737.6432500000001
$ws.Range("I5").Value = 0.8795180927887045
$ws.Range("J5").Value = 0.8795180927887046
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.3006476666666667
$ws.Range("N5").Value = 0.9019430000000001
$ws.Range("O5").Value = 0.2555551336960822
$ws.Range("P5").Value = 0.2555551336960821
$ws.Range("Q5").Value = 73.9235739816389
$ws.Range("R5").Value = 665.3121658347501
$ws.Range("S5").Value = 0.2247653637907406
$ws.Range("T5").Value = 0.2247653637907406

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ace"
$ws.Range("C6").Value = "Bdkrb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 245.8810833333334
$ws.Range("H6").Value = 737.6432500000001
$ws.Range("I6").Value = 0.8795180927887045
$ws.Range("J6").Value = 0.8795180927887046
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8000553333333333
$ws.Range("N6").Value = 2.400166
$ws.Range("O6").Value = 0.6800593197383766
$ws.Range("P6").Value = 0.6800593197383766
$ws.Range("Q6").Value = 196.7184720866111
$ws.Range("R6").Value = 1770.4662487795
$ws.Range("S6").Value = 0.5981244758794808
$ws.Range("T6").Value = 0.5981244758794809

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ace"
$ws.Range("C7").Value = "Bdkrb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 245.8810833333334
$ws.Range("H7").Value = 737.6432500000001
$ws.Range("I7").Value = 0.8795180927887045
$ws.Range("J7").Value = 0.8795180927887046
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.07574633333333333
$ws.Range("N7").Value = 0.227239
$ws.Range("O7").Value = 0.06438554656554128
$ws.Range("P7").Value = 0.06438554656554128
$ws.Range("Q7").Value = 18.62459049852778
$ws.Range("R7").Value = 167.62131448675
$ws.Range("S7").Value = 0.05662825311848319
$ws.Range("T7").Value = 0.0566282531184832

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ace"
$ws.Range("C8").Value = "Bdkrb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.178484
$ws.Range("H8").Value = 3.535452
$ws.Range("I8").Value = 0.004215444254639368
$ws.Range("J8").Value = 0.004215444254639368
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.3006476666666667
$ws.Range("N8").Value = 0.9019430000000001
$ws.Range("O8").Value = 0.2555551336960822
$ws.Range("P8").Value = 0.2555551336960821
$ws.Range("Q8").Value = 0.3543084648040001
$ws.Range("R8").Value = 3.188776183236
$ws.Range("S8").Value = 0.001077278420082745
$ws.Range("T8").Value = 0.001077278420082745

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ace"
$ws.Range("C9").Value = "Bdkrb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.178484
$ws.Range("H9").Value = 3.535452
$ws.Range("I9").Value = 0.004215444254639368
$ws.Range("J9").Value = 0.004215444254639368
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8000553333333333
$ws.Range("N9").Value = 2.400166
$ws.Range("O9").Value = 0.6800593197383766
$ws.Range("P9").Value = 0.6800593197383766
$ws.Range("Q9").Value = 0.9428524094480001
$ws.Range("R9").Value = 8.485671685032001
$ws.Range("S9").Value = 0.002866752152205097
$ws.Range("T9").Value = 0.002866752152205097

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ace"
$ws.Range("C10").Value = "Bdkrb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.178484
$ws.Range("H10").Value = 3.535452
$ws.Range("I10").Value = 0.004215444254639368
$ws.Range("J10").Value = 0.004215444254639368
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.07574633333333333
$ws.Range("N10").Value = 0.227239
$ws.Range("O10").Value = 0.06438554656554128
$ws.Range("P10").Value = 0.06438554656554128
$ws.Range("Q10").Value = 0.089265841892
$ws.Range("R10").Value = 0.803392577028
$ws.Range("S10").Value = 0.0002714136823515265
$ws.Range("T10").Value = 0.0002714136823515265

